$d = $word.ActiveDocument

$replacements = @(
    @("481×4=1924", "447×3=1341"),
    @("656×4=2624", "920×8=7360"),
    @("124×6=744", "917×2=1834"),
    @("497×4=1988", "270×8=2160"),
    @("349×9=3141", "788×5=3940"),
    @("851×5=4255", "559×5=2795"),
    @("868×5=4340", "323×8=2584"),
    @("723×2=1446", "741×8=5928"),
    @("449×6=2694", "230×6=1380"),
    @("964×3=2892", "249×4=996"),
    @("867×4=3468", "985×5=4925"),
    @("418×7=2926", "737×4=2948"),
    @("146×4=584", "381×5=1905"),
    @("411×3=1233", "889×6=5334"),
    @("973×9=8757", "272×6=1632"),
    @("594×3=1782", "790×7=5530"),
    @("735×2=1470", "319×4=1276"),
    @("518×2=1036", "595×7=4165"),
    @("714×3=2142", "494×9=4446"),
    @("290×7=2030", "529×8=4232"),
    @("137×9=1233", "799×7=5593"),
    @("436×5=2180", "504×4=2016"),
    @("296×5=1480", "455×5=2275"),
    @("658×3=1974", "315×3=945"),
    @("440×3=1320", "397×2=794")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Host "Done applying replacements"
